$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the report for week 3 (row 4): results and assignment notes
$ws.Range("F4").Value = "Thọ: Giao diện phần app.`nCông: Giao diện phần web API."
$ws.Range("E4").Value = "Hoàn thành xong giao diện."

# Match the wrap-text style used by the similar note cell F3
$ws.Range("F4").WrapText = $true

# Reflect the active selection recorded in the saved workbook
$ws.Range("E4").Select()
